# Auto-generated script to apply scheduled market-data update to Gilgamesh_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15 (Leve Item ID 44146)
$ws.Range("H15").Value = 3221.1785
$ws.Range("I15").Value = 3221.1785
$ws.Range("K15").Value = 9663.5355
$ws.Range("M15").Value = -9494.5355
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 11906532
$ws.Range("J17").Value = 11906532
$ws.Range("L17").Value = 35719596
$ws.Range("N17").Value = -35719932
# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 653.7646999999999
$ws.Range("J53").Value = 70
$ws.Range("L53").Value = 70
$ws.Range("N53").Value = -1344
# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 4060.2666
$ws.Range("I113").Value = 3652
$ws.Range("J113").Value = 4123.077
$ws.Range("K113").Value = 3652
$ws.Range("L113").Value = 4123.077
$ws.Range("M113").Value = -398
$ws.Range("N113").Value = -10631.077
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 10629.5
$ws.Range("I132").Value = 10629.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 31888.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -29358.5
$ws.Range("N132").ClearContents()
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 5338.9614
$ws.Range("I137").Value = 1793.2858
$ws.Range("J137").Value = 9475.583000000001
$ws.Range("K137").Value = 5379.857400000001
$ws.Range("L137").Value = 28426.749
$ws.Range("M137").Value = -2829.857400000001
$ws.Range("N137").Value = -33526.749

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3594.7432
$ws.Range("I32").Value = 2843.0286
$ws.Range("K32").Value = 2843.0286
$ws.Range("M32").Value = -2556.0286
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 29168.21
$ws.Range("J45").Value = 7693.4287
$ws.Range("L45").Value = 7693.4287
$ws.Range("N45").Value = -8447.4287
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 6002.273
$ws.Range("I61").Value = 2536.6667
$ws.Range("K61").Value = 2536.6667
$ws.Range("M61").Value = -2324.6667
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 5598.857
$ws.Range("I102").Value = 5354.8125
$ws.Range("K102").Value = 5354.8125
$ws.Range("M102").Value = -3732.8125
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2186.976
$ws.Range("I132").Value = 1481.9667
$ws.Range("J132").Value = 3949.5
$ws.Range("K132").Value = 4445.9001
$ws.Range("L132").Value = 11848.5
$ws.Range("M132").Value = -1915.9001
$ws.Range("N132").Value = -16908.5
# Row 135 (Leve Item ID 42016)
$ws.Range("H135").Value = 117982.664
$ws.Range("J135").Value = 117982.664
$ws.Range("L135").Value = 117982.664
$ws.Range("N135").Value = -128122.664
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 6002.273
$ws.Range("I136").Value = 2536.6667
$ws.Range("K136").Value = 7610.000100000001
$ws.Range("M136").Value = -5060.000100000001
# Row 139 (Leve Item ID 42321)
$ws.Range("H139").Value = 75748.336
$ws.Range("J139").Value = 75748.336
$ws.Range("L139").Value = 75748.336
$ws.Range("N139").Value = -86028.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 69 (Leve Item ID 10617)
$ws.Range("H69").Value = 19295
$ws.Range("J69").Value = 19295
$ws.Range("L69").Value = 19295
$ws.Range("N69").Value = -20917
# Row 72 (Leve Item ID 10617)
$ws.Range("H72").Value = 19295
$ws.Range("J72").Value = 19295
$ws.Range("L72").Value = 57885
$ws.Range("N72").Value = -65997
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 52631850
$ws.Range("I94").Value = 55555764
$ws.Range("J94").Value = 1399
$ws.Range("K94").Value = 55555764
$ws.Range("L94").Value = 1399
$ws.Range("M94").Value = -55555313
$ws.Range("N94").Value = -2301
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 81121.38
$ws.Range("I99").Value = 103557.8
$ws.Range("K99").Value = 103557.8
$ws.Range("M99").Value = -102059.8
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 13002755
$ws.Range("I105").Value = 771780.9399999999
$ws.Range("K105").Value = 771780.9399999999
$ws.Range("M105").Value = -770033.9399999999
# Row 122 (Leve Item ID 34096)
$ws.Range("H122").Value = 89997.5
$ws.Range("J122").Value = 89997.5
$ws.Range("L122").Value = 89997.5
$ws.Range("N122").Value = -99797.5
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2282.6191
$ws.Range("I134").Value = 1745.6666
$ws.Range("J134").Value = 3625
$ws.Range("K134").Value = 5236.9998
$ws.Range("L134").Value = 10875
$ws.Range("M134").Value = -2701.9998
$ws.Range("N134").Value = -15945

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2313.4
$ws.Range("I16").Value = 2287.2144
$ws.Range("K16").Value = 2287.2144
$ws.Range("M16").Value = -2000.2144
# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 10349.667
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2448
$ws.Range("I105").Value = 745
$ws.Range("J105").Value = 3177.8572
$ws.Range("K105").Value = 745
$ws.Range("L105").Value = 3177.8572
$ws.Range("M105").Value = 1002
$ws.Range("N105").Value = -6671.8572
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 3572097.8
$ws.Range("I107").Value = 4546124.5
$ws.Range("J107").Value = 666.3333
$ws.Range("K107").Value = 4546124.5
$ws.Range("L107").Value = 666.3333
$ws.Range("M107").Value = -4544204.5
$ws.Range("N107").Value = -4506.3333
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2313.4
$ws.Range("I113").Value = 2287.2144
$ws.Range("K113").Value = 2287.2144
$ws.Range("M113").Value = -117.2143999999998
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 6546.3335
$ws.Range("I134").Value = 6455.6
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 19366.8
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -16831.8
$ws.Range("N134").Value = -26070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3 (Leve Item ID 44094)
$ws.Range("H3").Value = 3523.375
$ws.Range("I3").Value = 2861.1428
$ws.Range("K3").Value = 8583.428400000001
$ws.Range("M3").Value = -8471.428400000001
# Row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 328.3
$ws.Range("I11").Value = 247.33333
$ws.Range("J11").Value = 449.75
$ws.Range("K11").Value = 741.99999
$ws.Range("L11").Value = 1349.25
$ws.Range("M11").Value = -601.99999
$ws.Range("N11").Value = -1629.25
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 7616.9414
$ws.Range("J113").Value = 9142.857
$ws.Range("L113").Value = 27428.571
$ws.Range("N113").Value = -31768.571
# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 3040.25
$ws.Range("I138").Value = 1795.7693
$ws.Range("K138").Value = 5387.3079
$ws.Range("M138").Value = -247.3078999999998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 111 (Leve Item ID 25853)
$ws.Range("H111").Value = 49949.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 49949.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 49949.5
$ws.Range("N111").Value = -56083.5
$ws.Range("M111").ClearContents()
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2763.2
$ws.Range("I132").Value = 2282.2222
$ws.Range("K132").Value = 6846.6666
$ws.Range("M132").Value = -4316.6666
# Row 141 (Leve Item ID 42504)
$ws.Range("H141").Value = 57170
$ws.Range("J141").Value = 57170
$ws.Range("L141").Value = 57170
$ws.Range("N141").Value = -67530

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 418.84616
$ws.Range("I22").Value = 448.72726
$ws.Range("K22").Value = 448.72726
$ws.Range("M22").Value = -153.72726
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 418.84616
$ws.Range("I27").Value = 448.72726
$ws.Range("K27").Value = 448.72726
$ws.Range("M27").Value = -341.72726
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 27659.28
$ws.Range("I40").Value = 29969.82
$ws.Range("K40").Value = 29969.82
$ws.Range("M40").Value = -29833.82
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3784.1428
$ws.Range("I68").Value = 3997.8
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 3997.8
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -3248.8
$ws.Range("N68").Value = -4748
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3784.1428
$ws.Range("I71").Value = 3997.8
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 19989
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -16245
$ws.Range("N71").Value = -23738
# Row 94 (Leve Item ID 18067)
$ws.Range("H94").Value = 1000000
$ws.Range("J94").Value = 1000000
$ws.Range("L94").Value = 1000000
$ws.Range("N94").Value = -1001352
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2122.8096
$ws.Range("I122").Value = 2071.5789
$ws.Range("K122").Value = 6214.736699999999
$ws.Range("M122").Value = -3764.736699999999
# Row 135 (Leve Item ID 42036)
$ws.Range("H135").Value = 109995.5
$ws.Range("J135").Value = 109995.5
$ws.Range("L135").Value = 109995.5
$ws.Range("N135").Value = -120135.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 8248.5
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 8248.5
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 624.26086
$ws.Range("I113").Value = 591.86664
$ws.Range("J113").Value = 685
$ws.Range("K113").Value = 1775.59992
$ws.Range("L113").Value = 2055
$ws.Range("M113").Value = 394.4000800000001
$ws.Range("N113").Value = -6395
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 8623665
$ws.Range("I122").Value = 2840.4783
$ws.Range("K122").Value = 8521.4349
$ws.Range("M122").Value = -6071.4349
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2284.25
$ws.Range("I126").Value = 2099.1765
$ws.Range("K126").Value = 6297.529500000001
$ws.Range("M126").Value = -3827.529500000001
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 11114477
$ws.Range("I132").Value = 13892230
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 41676690
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -41674160
$ws.Range("N132").Value = -15458.9999
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 27029714
$ws.Range("I136").Value = 31251206
$ws.Range("K136").Value = 93753618
$ws.Range("M136").Value = -93751068
# Row 140 (Leve Item ID 42506)
$ws.Range("H140").Value = 109562.71
$ws.Range("J140").Value = 109562.71
$ws.Range("L140").Value = 109562.71
$ws.Range("N140").Value = -119922.71
# Row 141 (Leve Item ID 42505)
$ws.Range("H141").Value = 70141.57000000001
$ws.Range("J141").Value = 70141.57000000001
$ws.Range("L141").Value = 70141.57000000001
$ws.Range("N141").Value = -80501.57000000001

